$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# 1) Fix the double-logged / border-case text in row 71 (F71): add "PICK UP 6 MATS."
$ws.Range("F71").Value = "Pick up PC and Projector cart. Projector cart has small speaker on it. Pick up portable screen also. PICK UP 6 MATS. Key for room in Founders 164 storeroom. Return equipment to Vanier 040 storeroom."
$ws.Rows(71).RowHeight = 43.15

# 2) Append new log rows (73-76) for the next day (Monday) of entries.
#    Row 73 is a day-separator row, matching the style of row 69 ("FRIDAY").
$ws.Rows(69).Copy()
$ws.Rows(73).PasteSpecial(-4122)
$ws.Range("C73").Value = "MONDAY"

$ws.Range("A74").Value = "Other"
$ws.Range("B74").Value = 42604
$ws.Range("C74").Value = "1530"
$ws.Range("D74").Value = "MC"
$ws.Range("E74").Value = "157A"
$ws.Range("F74").Value = ""

$ws.Range("A75").Value = "Other"
$ws.Range("B75").Value = 42604
$ws.Range("C75").Value = "1530"
$ws.Range("D75").Value = "WC"
$ws.Range("E75").Value = "117"
$ws.Range("F75").Value = "Winters classroom key B15"

$ws.Range("A76").Value = "Other"
$ws.Range("B76").Value = 42604
$ws.Range("C76").Value = "1600"
$ws.Range("D76").Value = "WC"
$ws.Range("E76").Value = "283B"
$ws.Range("F76").Value = "No need to go here. Room is managed by Winters college."

# Match style / number format of rows 74-76 to the existing equivalent rows (70/71/72 pattern)
$ws.Range("A74").Style = $ws.Range("A71").Style
$ws.Range("B74").NumberFormat = $ws.Range("B71").NumberFormat
$ws.Range("C74").Style = $ws.Range("C8").Style
$ws.Range("D74").Style = $ws.Range("D8").Style
$ws.Range("E74").Style = $ws.Range("E19").Style
$ws.Range("F74").Style = $ws.Range("F19").Style

$ws.Range("A75").Style = $ws.Range("A71").Style
$ws.Range("B75").NumberFormat = $ws.Range("B71").NumberFormat
$ws.Range("C75").Style = $ws.Range("C71").Style
$ws.Range("D75").Style = $ws.Range("D71").Style
$ws.Range("E75").Style = $ws.Range("E71").Style
$ws.Range("F75").Style = $ws.Range("F71").Style

$ws.Range("A76").Style = $ws.Range("A71").Style
$ws.Range("B76").NumberFormat = $ws.Range("B71").NumberFormat
$ws.Range("C76").Style = $ws.Range("C71").Style
$ws.Range("D76").Style = $ws.Range("D71").Style
$ws.Range("E76").Style = $ws.Range("E71").Style
$ws.Range("F76").Style = $ws.Range("F71").Style

# 3) Update the frozen pane / active selection to reflect the scrolled-down view
$win = $excel.ActiveWindow
$win.ScrollRow = 60
$ws.Range("F82").Select()
